$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (values copied from source row's
# original D/J/K/L/M/P tuple into destination row), derived from the diff.
$rowMap = @{
    2  = 9
    3  = 17
    4  = 2
    5  = 4
    6  = 15
    7  = 14
    8  = 16
    9  = 10
    10 = 6
    12 = 7
    13 = 18
    14 = 12
    15 = 5
    16 = 13
    17 = 8
    18 = 3
}

# Capture the original values for the columns that move between rows before
# overwriting anything, so the permutation is applied consistently.
$cols = @(4, 10, 11, 12, 13, 16)  # D, J, K, L, M, P
$original = @{}
foreach ($r in $rowMap.Keys) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $vals
}

foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $original[$src][$c]
    }
}
